$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume/Number and date range) ---
# A8: "Volume 30   Number  16" -> "Volume 30   Number  17"
$ws.Range("A8").Value = "Volume 30   Number  17"
# C9: "Report Covering the Week  4/17/2023  Through  4/23/2023" -> "...4/24/2023  Through  4/30/2023"
$ws.Range("C9").Value = "Report Covering the Week  4/24/2023  Through  4/30/2023"

# --- Weekly crime statistics table updates (rows 14-30) ---
$ws.Range("C14").Value = 3
$ws.Range("D14").Value = 11
$ws.Range("E14").Value = -72.727272727272
$ws.Range("F14").Value = 29
$ws.Range("G14").Value = 33
$ws.Range("H14").Value = -12.121212121212
$ws.Range("I14").Value = 125
$ws.Range("J14").Value = 136
$ws.Range("K14").Value = -8.088235294117
$ws.Range("L14").Value = -16.666666666666
$ws.Range("M14").Value = -17.218543046357
$ws.Range("N14").Value = -79.740680713128
$ws.Range("C15").Value = 35
$ws.Range("D15").Value = 18
$ws.Range("E15").Value = 94.444444444444
$ws.Range("F15").Value = 116
$ws.Range("G15").Value = 101
$ws.Range("H15").Value = 14.851485148514
$ws.Range("I15").Value = 494
$ws.Range("J15").Value = 511
$ws.Range("K15").Value = -3.326810176125
$ws.Range("L15").Value = 8.333333333333
$ws.Range("M15").Value = 24.120603015075
$ws.Range("N15").Value = -50.451354062186
$ws.Range("C16").Value = 284
$ws.Range("D16").Value = 276
$ws.Range("E16").Value = 2.898550724637
$ws.Range("F16").Value = 1136
$ws.Range("G16").Value = 1182
$ws.Range("H16").Value = -3.891708967851
$ws.Range("I16").Value = 4991
$ws.Range("J16").Value = 5094
$ws.Range("K16").Value = -2.021986650961
$ws.Range("L16").Value = 42.234254773439
$ws.Range("M16").Value = -13.184901722038
$ws.Range("N16").Value = -81.736012002781
$ws.Range("C17").Value = 446
$ws.Range("D17").Value = 481
$ws.Range("E17").Value = -7.276507276507
$ws.Range("F17").Value = 2072
$ws.Range("G17").Value = 1973
$ws.Range("H17").Value = 5.01773948302
$ws.Range("I17").Value = 8352
$ws.Range("J17").Value = 7682
$ws.Range("K17").Value = 8.721687060661
$ws.Range("L17").Value = 32.424290470905
$ws.Range("M17").Value = 62.237762237762
$ws.Range("N17").Value = -31.072047536518
$ws.Range("C18").Value = 254
$ws.Range("D18").Value = 298
$ws.Range("E18").Value = -14.76510067114
$ws.Range("F18").Value = 1022
$ws.Range("G18").Value = 1165
$ws.Range("H18").Value = -12.274678111588
$ws.Range("I18").Value = 4673
$ws.Range("J18").Value = 5041
$ws.Range("K18").Value = -7.300138861337
$ws.Range("L18").Value = 24.546908315565
$ws.Range("M18").Value = -17.321302193913
$ws.Range("N18").Value = -85.523992441374
$ws.Range("C19").Value = 922
$ws.Range("D19").Value = 966
$ws.Range("E19").Value = -4.55486542443
$ws.Range("F19").Value = 3764
$ws.Range("G19").Value = 3668
$ws.Range("H19").Value = 2.617230098146
$ws.Range("I19").Value = 15706
$ws.Range("J19").Value = 15934
$ws.Range("K19").Value = -1.430902472699
$ws.Range("L19").Value = 52.692980750534
$ws.Range("M19").Value = 38.342288381925
$ws.Range("N19").Value = -39.326276751912
$ws.Range("C20").Value = 304
$ws.Range("D20").Value = 242
$ws.Range("E20").Value = 25.619834710743
$ws.Range("F20").Value = 1231
$ws.Range("G20").Value = 877
$ws.Range("H20").Value = 40.36488027366
$ws.Range("I20").Value = 4793
$ws.Range("J20").Value = 4208
$ws.Range("K20").Value = 13.902091254752
$ws.Range("L20").Value = 88.108320251177
$ws.Range("M20").Value = 50.392218387198
$ws.Range("N20").Value = -87.149100463844
$ws.Range("C21").Value = 2248
$ws.Range("D21").Value = 2292
$ws.Range("E21").Value = -1.919720767888
$ws.Range("F21").Value = 9370
$ws.Range("G21").Value = 8999
$ws.Range("H21").Value = 4.12268029781
$ws.Range("I21").Value = 39134
$ws.Range("J21").Value = 38606
$ws.Range("K21").Value = 1.367663057555
$ws.Range("L21").Value = 44.897808056872
$ws.Range("M21").Value = 23.693027372147
$ws.Range("N21").Value = -71.335022926707
$ws.Range("C22").Value = 34
$ws.Range("D22").Value = 47
$ws.Range("E22").Value = -27.659574468085
$ws.Range("F22").Value = 157
$ws.Range("G22").Value = 170
$ws.Range("H22").Value = -7.647058823529
$ws.Range("I22").Value = 696
$ws.Range("J22").Value = 757
$ws.Range("K22").Value = -8.058124174372
$ws.Range("L22").Value = 53.982300884955
$ws.Range("M22").Value = 3.111111111111
$ws.Range("C23").Value = 101
$ws.Range("D23").Value = 110
$ws.Range("E23").Value = -8.181818181818
$ws.Range("F23").Value = 462
$ws.Range("G23").Value = 426
$ws.Range("H23").Value = 8.450704225352
$ws.Range("I23").Value = 1972
$ws.Range("J23").Value = 1811
$ws.Range("K23").Value = 8.890115958034
$ws.Range("L23").Value = 20.611620795107
$ws.Range("M23").Value = 66.413502109704
$ws.Range("C24").Value = 2133
$ws.Range("D24").Value = 2242
$ws.Range("E24").Value = -4.86173059768
$ws.Range("F24").Value = 8203
$ws.Range("G24").Value = 8767
$ws.Range("H24").Value = -6.433215467092
$ws.Range("I24").Value = 34738
$ws.Range("J24").Value = 34845
$ws.Range("K24").Value = -0.307074185679
$ws.Range("L24").Value = 42.130027412953
$ws.Range("M24").Value = 42.066088663504
$ws.Range("C25").Value = 838
$ws.Range("D25").Value = 839
$ws.Range("E25").Value = -0.119189511323
$ws.Range("F25").Value = 3315
$ws.Range("G25").Value = 3210
$ws.Range("H25").Value = 3.271028037383
$ws.Range("I25").Value = 13458
$ws.Range("J25").Value = 12786
$ws.Range("K25").Value = 5.255748474894
$ws.Range("L25").Value = 36.643314042034
$ws.Range("M25").Value = -4.172600398746
$ws.Range("C26").Value = 52
$ws.Range("D26").Value = 38
$ws.Range("E26").Value = 36.842105263157
$ws.Range("F26").Value = 183
$ws.Range("G26").Value = 176
$ws.Range("H26").Value = 3.977272727272
$ws.Range("I26").Value = 787
$ws.Range("J26").Value = 833
$ws.Range("K26").Value = -5.522208883553
$ws.Range("L26").Value = 2.473958333333
$ws.Range("C27").Value = 105
$ws.Range("D27").Value = 111
$ws.Range("E27").Value = -5.405405405405
$ws.Range("F27").Value = 416
$ws.Range("G27").Value = 382
$ws.Range("H27").Value = 8.900523560209
$ws.Range("I27").Value = 1635
$ws.Range("J27").Value = 1541
$ws.Range("K27").Value = 6.099935107073
$ws.Range("L27").Value = 18.135838150289
$ws.Range("D28").Value = 25
$ws.Range("E28").Value = -48
$ws.Range("F28").Value = 73
$ws.Range("G28").Value = 132
$ws.Range("H28").Value = -44.696969696969
$ws.Range("I28").Value = 349
$ws.Range("J28").Value = 461
$ws.Range("K28").Value = -24.295010845987
$ws.Range("L28").Value = -21.218961625282
$ws.Range("M28").Value = -27.291666666666
$ws.Range("N28").Value = -80.349099099099
$ws.Range("C29").Value = 13
$ws.Range("D29").Value = 22
$ws.Range("E29").Value = -40.90909090909
$ws.Range("F29").Value = 65
$ws.Range("G29").Value = 101
$ws.Range("H29").Value = -35.643564356435
$ws.Range("I29").Value = 294
$ws.Range("J29").Value = 394
$ws.Range("K29").Value = -25.380710659898
$ws.Range("L29").Value = -26.5
$ws.Range("M29").Value = -25.569620253164
$ws.Range("N29").Value = -81.851851851851
$ws.Range("C30").Value = 10
$ws.Range("D30").Value = 15
$ws.Range("E30").Value = -33.333333333333
$ws.Range("F30").Value = 39
$ws.Range("G30").Value = 38
$ws.Range("H30").Value = 2.631578947368
$ws.Range("I30").Value = 155
$ws.Range("J30").Value = 230
$ws.Range("K30").Value = -32.608695652173
$ws.Range("L30").Value = 1.307189542483
